$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Dated 10/30/2017" -> "Last updated 10/29/2018", and move the _GoBack
#    bookmark from the end of the "...baseline assessments." paragraph to
#    right after this run (still zero-width, sitting after the run, inside
#    the same paragraph).
# ---------------------------------------------------------------------------

# Replace the text first, appending a temporary marker character "X" so that
# the bookmark-insertion point below does not land exactly on the boundary
# of the freshly-edited run (landing exactly there is mishandled).
$d.Content.Find.Execute("Dated 10/30/2017", $false, $false, $false, $false, $false, $true, 1, $false, "Last updated 10/29/2018X", 2)

# Locate the "X" marker and the position right before it.
$rMarker = $d.Content
$rMarker.Find.Execute("2018X", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markerStart = $rMarker.End - 1

# Add (== relocate, since the name already exists elsewhere in the doc) the
# _GoBack bookmark as a zero-width point right before the "X" marker.
$rPoint = $d.Range($markerStart, $markerStart)
$d.Bookmarks.Add("_GoBack", $rPoint)

# Remove the temporary marker character.
$rX = $d.Range($markerStart, $markerStart + 1)
$rX.Delete()

# ---------------------------------------------------------------------------
# 2) Restructure "... between the 2 groups usual care ..." so that the
#    grammar-check proofErr markers wrap "2" instead of "groups" (net visible
#    text is unchanged).
# ---------------------------------------------------------------------------

# Find the existing "groups" run, which is already wrapped by the
# gramStart/gramEnd proofErr markers, and swap its text for "2" so the
# existing proofErr wrapping simply carries over onto the new content.
$rGroups = $d.Content
$rGroups.Find.Execute("the 2 groups usual care and the HELP clinic", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
# rGroups now spans "the 2 groups usual care and the HELP clinic"; narrow it
# down precisely using the known sub-offsets.
$oldTwoStart = $rGroups.Start + 4            # position of "2"
$oldGroupsStart = $oldTwoStart + 2           # position of "groups"
$oldGroupsEnd = $oldGroupsStart + 6          # end of "groups"

$rWord = $d.Range($oldGroupsStart, $oldGroupsEnd)
$rWord.Text = "2"

$rOldTwo = $d.Range($oldTwoStart, $oldTwoStart + 2)
$rOldTwo.Delete()

$rInsert = $d.Range($oldTwoStart, $oldTwoStart)
$rInsert.InsertBefore("groups ")

# ---------------------------------------------------------------------------
# 3) Remove the now-redundant *old* _GoBack bookmark location (it has already
#    been relocated by step 1's Bookmarks.Add call, so nothing further is
#    required here, but guard in case only one survives).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 4) Merge the three runs "(3", ")=", "1.066" into a single run "(3)=1.066".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("(3)=1.066", $false, $false, $false, $false, $false, $true, 1, $false, "(3)=1.066", 2)

# ---------------------------------------------------------------------------
# 5) Merge the three runs that spell out the "239 subjects have missing
#    data..." sentence opening into a single run.
# ---------------------------------------------------------------------------
$quote = [char]0x201C
$mergedText = "239 subjects have missing data for " + $quote + "Number of times in detox "
$d.Content.Find.Execute($mergedText, $false, $false, $false, $false, $false, $true, 1, $false, $mergedText, 2)
